# "Quick save before branch": duplicate sheet "56" into a new trailing
# sheet "57" (same layout/data/styles), bump its filenumber cell (B1)
# from 56 to 57, and leave it as the active/selected sheet+tab.

$wb = $excel.ActiveWorkbook

$src = $wb.Worksheets.Item("56")

# Copy the "56" sheet, placing the new copy immediately after it.
# Excel names the duplicate "56 (2)" by default; it also becomes the
# active sheet/tab, matching the workbook's new activeTab.
$src.Copy($null, $src)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "57"

# Update the duplicated sheet's filenumber value to match its new name.
$newSheet.Range("B1").Value = 57
